$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 493, shifting existing rows 493:523 down to 494:524.
$ws.Rows.Item(493).Insert()

# Populate the newly inserted row 493 with the new record's data.
$ws.Cells.Item(493, 1).Value = 3
$ws.Cells.Item(493, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(493, 3).Value = "Coquimbo"
$ws.Cells.Item(493, 4).Value = 44931
$ws.Cells.Item(493, 5).Value = 5
$ws.Cells.Item(493, 6).Value = 100112017
$ws.Cells.Item(493, 7).Value = "Apio"
$ws.Cells.Item(493, 8).Value = "Americana (o)"
$ws.Cells.Item(493, 9).Value = "Primera"
$ws.Cells.Item(493, 10).Value = 165
$ws.Cells.Item(493, 11).Value = 12000
$ws.Cells.Item(493, 12).Value = 13000
$ws.Cells.Item(493, 13).Value = 12515
$ws.Cells.Item(493, 14).Value = "$/docena de matas"
$ws.Cells.Item(493, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(493, 16).Value = 2086
$ws.Cells.Item(493, 17).Value = 6
$ws.Cells.Item(493, 18).Value = "Hortaliza"
